$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 90927  # H11
$ws.Cells.Item(11, 9).Value = 90927  # I11
$ws.Cells.Item(11, 11).Value = 90927  # K11
$ws.Cells.Item(11, 13).Value = -90787  # M11
$ws.Cells.Item(43, 8).Value = 4166.6665  # H43
$ws.Cells.Item(43, 9).Value = 3375  # I43
$ws.Cells.Item(43, 10).Value = 5750  # J43
$ws.Cells.Item(43, 11).Value = 3375  # K43
$ws.Cells.Item(43, 12).Value = 5750  # L43
$ws.Cells.Item(43, 13).Value = -3306  # M43
$ws.Cells.Item(43, 14).Value = -5888  # N43
$ws.Cells.Item(68, 8).Value = 36200  # H68
$ws.Cells.Item(68, 10).Value = 36200  # J68
$ws.Cells.Item(68, 12).Value = 36200  # L68
$ws.Cells.Item(68, 14).Value = -37698  # N68
$ws.Cells.Item(71, 8).Value = 36200  # H71
$ws.Cells.Item(71, 10).Value = 36200  # J71
$ws.Cells.Item(71, 12).Value = 108600  # L71
$ws.Cells.Item(71, 14).Value = -116088  # N71
$ws.Cells.Item(103, 8).Value = 1434.6471  # H103
$ws.Cells.Item(103, 9).Value = 1989  # I103
$ws.Cells.Item(103, 10).Value = 1132.2727  # J103
$ws.Cells.Item(103, 11).Value = 5967  # K103
$ws.Cells.Item(103, 12).Value = 3396.8181  # L103
$ws.Cells.Item(103, 13).Value = -5381  # M103
$ws.Cells.Item(103, 14).Value = -4568.8181  # N103
$ws.Cells.Item(121, 8).Value = 649.875  # H121
$ws.Cells.Item(121, 9).Value = 0  # I121
$ws.Cells.Item(121, 10).Value = 649.875  # J121
$ws.Cells.Item(121, 11).Value = 0  # K121
$ws.Cells.Item(121, 12).Value = 1949.625  # L121
$ws.Cells.Item(121, 13).ClearContents()  # M121
$ws.Cells.Item(121, 14).Value = -5443.625  # N121
$ws.Cells.Item(131, 8).Value = 3807.6033  # H131
$ws.Cells.Item(131, 9).Value = 777.93335  # I131
$ws.Cells.Item(131, 10).Value = 4754.375  # J131
$ws.Cells.Item(131, 11).Value = 2333.80005  # K131
$ws.Cells.Item(131, 12).Value = 14263.125  # L131
$ws.Cells.Item(131, 13).Value = 2706.19995  # M131
$ws.Cells.Item(131, 14).Value = -24343.125  # N131
$ws.Cells.Item(132, 8).Value = 5573.8  # H132
$ws.Cells.Item(132, 9).Value = 5572.7593  # I132
$ws.Cells.Item(132, 11).Value = 16718.2779  # K132
$ws.Cells.Item(132, 13).Value = -14188.2779  # M132
$ws.Cells.Item(135, 8).Value = 3441.4055  # H135
$ws.Cells.Item(135, 9).Value = 926.3  # I135
$ws.Cells.Item(135, 10).Value = 6400.353  # J135
$ws.Cells.Item(135, 11).Value = 8336.699999999999  # K135
$ws.Cells.Item(135, 12).Value = 57603.177  # L135
$ws.Cells.Item(135, 13).Value = -5801.699999999999  # M135
$ws.Cells.Item(135, 14).Value = -62673.177  # N135
$ws.Cells.Item(137, 8).Value = 1380.975  # H137
$ws.Cells.Item(137, 9).Value = 1035.4  # I137
$ws.Cells.Item(137, 10).Value = 3800  # J137
$ws.Cells.Item(137, 11).Value = 3106.2  # K137
$ws.Cells.Item(137, 12).Value = 11400  # L137
$ws.Cells.Item(137, 13).Value = -556.2000000000003  # M137
$ws.Cells.Item(137, 14).Value = -16500  # N137
$ws.Cells.Item(138, 8).Value = 3070.524  # H138
$ws.Cells.Item(138, 9).Value = 1696.5358  # I138
$ws.Cells.Item(138, 10).Value = 5818.5  # J138
$ws.Cells.Item(138, 11).Value = 5089.607400000001  # K138
$ws.Cells.Item(138, 12).Value = 17455.5  # L138
$ws.Cells.Item(138, 13).Value = 50.39259999999922  # M138
$ws.Cells.Item(138, 14).Value = -27735.5  # N138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 29400  # H24
$ws.Cells.Item(24, 10).Value = 29400  # J24
$ws.Cells.Item(24, 12).Value = 29400  # L24
$ws.Cells.Item(24, 14).Value = -30148  # N24
$ws.Cells.Item(32, 8).Value = 27283.361  # H32
$ws.Cells.Item(32, 9).Value = 11015.586  # I32
$ws.Cells.Item(32, 10).Value = 74731.03999999999  # J32
$ws.Cells.Item(32, 11).Value = 11015.586  # K32
$ws.Cells.Item(32, 12).Value = 74731.03999999999  # L32
$ws.Cells.Item(32, 13).Value = -10728.586  # M32
$ws.Cells.Item(32, 14).Value = -75305.03999999999  # N32
$ws.Cells.Item(55, 8).Value = 14762.5  # H55
$ws.Cells.Item(55, 10).Value = 15585.714  # J55
$ws.Cells.Item(55, 12).Value = 15585.714  # L55
$ws.Cells.Item(55, 14).Value = -16215.714  # N55
$ws.Cells.Item(74, 8).Value = 1345.8  # H74
$ws.Cells.Item(74, 9).Value = 1334.6818  # I74
$ws.Cells.Item(74, 10).Value = 1376.375  # J74
$ws.Cells.Item(74, 11).Value = 1334.6818  # K74
$ws.Cells.Item(74, 12).Value = 1376.375  # L74
$ws.Cells.Item(74, 13).Value = -460.6818000000001  # M74
$ws.Cells.Item(74, 14).Value = -3124.375  # N74
$ws.Cells.Item(77, 8).Value = 1345.8  # H77
$ws.Cells.Item(77, 9).Value = 1334.6818  # I77
$ws.Cells.Item(77, 10).Value = 1376.375  # J77
$ws.Cells.Item(77, 11).Value = 6673.409000000001  # K77
$ws.Cells.Item(77, 12).Value = 6881.875  # L77
$ws.Cells.Item(77, 13).Value = -2305.409000000001  # M77
$ws.Cells.Item(77, 14).Value = -15617.875  # N77
$ws.Cells.Item(80, 8).Value = 25215.2  # H80
$ws.Cells.Item(80, 10).Value = 25215.2  # J80
$ws.Cells.Item(80, 12).Value = 25215.2  # L80
$ws.Cells.Item(80, 14).Value = -27211.2  # N80
$ws.Cells.Item(83, 8).Value = 25215.2  # H83
$ws.Cells.Item(83, 10).Value = 25215.2  # J83
$ws.Cells.Item(83, 12).Value = 75645.60000000001  # L83
$ws.Cells.Item(83, 14).Value = -85629.60000000001  # N83
$ws.Cells.Item(100, 8).Value = 29400  # H100
$ws.Cells.Item(100, 10).Value = 29400  # J100
$ws.Cells.Item(100, 12).Value = 29400  # L100
$ws.Cells.Item(100, 14).Value = -31564  # N100
$ws.Cells.Item(102, 8).Value = 69694.39999999999  # H102
$ws.Cells.Item(102, 9).Value = 102048  # I102
$ws.Cells.Item(102, 10).Value = 4987.2  # J102
$ws.Cells.Item(102, 11).Value = 102048  # K102
$ws.Cells.Item(102, 12).Value = 4987.2  # L102
$ws.Cells.Item(102, 13).Value = -100426  # M102
$ws.Cells.Item(102, 14).Value = -8231.200000000001  # N102

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 46765.477  # H20
$ws.Cells.Item(20, 9).Value = 56205.58  # I20
$ws.Cells.Item(20, 11).Value = 56205.58  # K20
$ws.Cells.Item(20, 13).Value = -55958.58  # M20

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 43833.37  # H31
$ws.Cells.Item(31, 9).Value = 1202.8572  # I31
$ws.Cells.Item(31, 10).Value = 72253.71000000001  # J31
$ws.Cells.Item(31, 11).Value = 1202.8572  # K31
$ws.Cells.Item(31, 12).Value = 72253.71000000001  # L31
$ws.Cells.Item(31, 13).Value = -907.8571999999999  # M31
$ws.Cells.Item(31, 14).Value = -72843.71000000001  # N31
$ws.Cells.Item(34, 8).Value = 43833.37  # H34
$ws.Cells.Item(34, 9).Value = 1202.8572  # I34
$ws.Cells.Item(34, 10).Value = 72253.71000000001  # J34
$ws.Cells.Item(34, 11).Value = 1202.8572  # K34
$ws.Cells.Item(34, 12).Value = 72253.71000000001  # L34
$ws.Cells.Item(34, 13).Value = -1000.8572  # M34
$ws.Cells.Item(34, 14).Value = -72657.71000000001  # N34
$ws.Cells.Item(45, 8).Value = 13749.5  # H45
$ws.Cells.Item(45, 9).Value = 12499  # I45
$ws.Cells.Item(45, 11).Value = 12499  # K45
$ws.Cells.Item(45, 13).Value = -11906  # M45
$ws.Cells.Item(58, 8).Value = 1491.561  # H58
$ws.Cells.Item(58, 10).Value = 2132.4  # J58
$ws.Cells.Item(58, 12).Value = 2132.4  # L58
$ws.Cells.Item(58, 14).Value = -2538.4  # N58
$ws.Cells.Item(94, 8).Value = 1279.6364  # H94
$ws.Cells.Item(94, 10).Value = 1302.5625  # J94
$ws.Cells.Item(94, 12).Value = 1302.5625  # L94
$ws.Cells.Item(94, 14).Value = -2204.5625  # N94
$ws.Cells.Item(131, 8).Value = 24635.408  # H131
$ws.Cells.Item(131, 10).Value = 24635.408  # J131
$ws.Cells.Item(131, 12).Value = 24635.408  # L131
$ws.Cells.Item(131, 14).Value = -34715.408  # N131
$ws.Cells.Item(132, 8).Value = 1932.7797  # H132
$ws.Cells.Item(132, 9).Value = 1831  # I132
$ws.Cells.Item(132, 10).Value = 2259.9285  # J132
$ws.Cells.Item(132, 11).Value = 5493  # K132
$ws.Cells.Item(132, 12).Value = 6779.7855  # L132
$ws.Cells.Item(132, 13).Value = -2963  # M132
$ws.Cells.Item(132, 14).Value = -11839.7855  # N132
$ws.Cells.Item(134, 8).Value = 914.9  # H134
$ws.Cells.Item(134, 9).Value = 545.0417  # I134
$ws.Cells.Item(134, 10).Value = 2394.3333  # J134
$ws.Cells.Item(134, 11).Value = 1635.1251  # K134
$ws.Cells.Item(134, 12).Value = 7182.999899999999  # L134
$ws.Cells.Item(134, 13).Value = 899.8749  # M134
$ws.Cells.Item(134, 14).Value = -12252.9999  # N134
$ws.Cells.Item(136, 8).Value = 1491.561  # H136
$ws.Cells.Item(136, 10).Value = 2132.4  # J136
$ws.Cells.Item(136, 12).Value = 6397.200000000001  # L136
$ws.Cells.Item(136, 14).Value = -11497.2  # N136

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1172.9  # H5
$ws.Cells.Item(5, 9).Value = 798.5599999999999  # I5
$ws.Cells.Item(5, 10).Value = 1440.2858  # J5
$ws.Cells.Item(5, 11).Value = 2395.68  # K5
$ws.Cells.Item(5, 12).Value = 4320.857400000001  # L5
$ws.Cells.Item(5, 13).Value = -2283.68  # M5
$ws.Cells.Item(5, 14).Value = -4544.857400000001  # N5
$ws.Cells.Item(107, 8).Value = 1166.5834  # H107
$ws.Cells.Item(107, 9).Value = 535.17645  # I107
$ws.Cells.Item(107, 10).Value = 2700  # J107
$ws.Cells.Item(107, 11).Value = 1605.52935  # K107
$ws.Cells.Item(107, 12).Value = 8100  # L107
$ws.Cells.Item(107, 13).Value = 314.4706499999998  # M107
$ws.Cells.Item(107, 14).Value = -11940  # N107
$ws.Cells.Item(113, 8).Value = 961.7778  # H113
$ws.Cells.Item(113, 9).Value = 1262  # I113
$ws.Cells.Item(113, 10).Value = 638.46155  # J113
$ws.Cells.Item(113, 11).Value = 3786  # K113
$ws.Cells.Item(113, 12).Value = 1915.38465  # L113
$ws.Cells.Item(113, 13).Value = -1616  # M113
$ws.Cells.Item(113, 14).Value = -6255.38465  # N113
$ws.Cells.Item(131, 8).Value = 822.33  # H131
$ws.Cells.Item(131, 9).Value = 563.1875  # I131
$ws.Cells.Item(131, 10).Value = 871.6905  # J131
$ws.Cells.Item(131, 11).Value = 1689.5625  # K131
$ws.Cells.Item(131, 12).Value = 2615.0715  # L131
$ws.Cells.Item(131, 13).Value = 3350.4375  # M131
$ws.Cells.Item(131, 14).Value = -12695.0715  # N131
$ws.Cells.Item(132, 8).Value = 501751.16  # H132
$ws.Cells.Item(132, 9).Value = 812.6667  # I132
$ws.Cells.Item(132, 11).Value = 7314.0003  # K132
$ws.Cells.Item(132, 13).Value = -4784.0003  # M132
$ws.Cells.Item(135, 8).Value = 1172.9  # H135
$ws.Cells.Item(135, 9).Value = 798.5599999999999  # I135
$ws.Cells.Item(135, 10).Value = 1440.2858  # J135
$ws.Cells.Item(135, 11).Value = 7187.039999999999  # K135
$ws.Cells.Item(135, 12).Value = 12962.5722  # L135
$ws.Cells.Item(135, 13).Value = -4652.039999999999  # M135
$ws.Cells.Item(135, 14).Value = -18032.5722  # N135
$ws.Cells.Item(136, 8).Value = 2073.6365  # H136
$ws.Cells.Item(136, 9).Value = 1883  # I136
$ws.Cells.Item(136, 11).Value = 5649  # K136
$ws.Cells.Item(136, 13).Value = -549  # M136

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(69, 8).Value = 100000  # H69
$ws.Cells.Item(69, 10).Value = 100000  # J69
$ws.Cells.Item(69, 12).Value = 100000  # L69
$ws.Cells.Item(69, 14).Value = -101498  # N69
$ws.Cells.Item(72, 8).Value = 100000  # H72
$ws.Cells.Item(72, 10).Value = 100000  # J72
$ws.Cells.Item(72, 12).Value = 300000  # L72
$ws.Cells.Item(72, 14).Value = -307488  # N72
$ws.Cells.Item(97, 8).Value = 100002216  # H97
$ws.Cells.Item(97, 9).Value = 125002500  # I97
$ws.Cells.Item(97, 10).Value = 1100  # J97
$ws.Cells.Item(97, 11).Value = 125002500  # K97
$ws.Cells.Item(97, 12).Value = 1100  # L97
$ws.Cells.Item(97, 13).Value = -125002004  # M97
$ws.Cells.Item(97, 14).Value = -2092  # N97
$ws.Cells.Item(132, 8).Value = 3532.8  # H132
$ws.Cells.Item(132, 9).Value = 3587.6667  # I132
$ws.Cells.Item(132, 10).Value = 3501.9375  # J132
$ws.Cells.Item(132, 11).Value = 10763.0001  # K132
$ws.Cells.Item(132, 12).Value = 10505.8125  # L132
$ws.Cells.Item(132, 13).Value = -8233.000100000001  # M132
$ws.Cells.Item(132, 14).Value = -15565.8125  # N132

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3600.2104  # H122
$ws.Cells.Item(122, 9).Value = 3517  # I122
$ws.Cells.Item(122, 10).Value = 3742.8572  # J122
$ws.Cells.Item(122, 11).Value = 10551  # K122
$ws.Cells.Item(122, 12).Value = 11228.5716  # L122
$ws.Cells.Item(122, 13).Value = -8101  # M122
$ws.Cells.Item(122, 14).Value = -16128.5716  # N122
$ws.Cells.Item(132, 8).Value = 3249.0286  # H132
$ws.Cells.Item(132, 9).Value = 3478.8262  # I132
$ws.Cells.Item(132, 10).Value = 2808.5833  # J132
$ws.Cells.Item(132, 11).Value = 10436.4786  # K132
$ws.Cells.Item(132, 12).Value = 8425.749899999999  # L132
$ws.Cells.Item(132, 13).Value = -7906.4786  # M132
$ws.Cells.Item(132, 14).Value = -13485.7499  # N132

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4673.449  # H132
$ws.Cells.Item(132, 9).Value = 3435.0908  # I132
$ws.Cells.Item(132, 10).Value = 9538.429  # J132
$ws.Cells.Item(132, 11).Value = 10305.2724  # K132
$ws.Cells.Item(132, 12).Value = 28615.287  # L132
$ws.Cells.Item(132, 13).Value = -7775.2724  # M132
$ws.Cells.Item(132, 14).Value = -33675.287  # N132
